$wb = $excel.ActiveWorkbook

# --- "package-private" sheet (sheet2.xml): add Reference/Components annotations ---
$pp = $wb.Worksheets.Item("package-private")

$pp.Range("D3").Value  = "Effectively immutable"
$pp.Range("D4").Value  = "Effectively immutable"
$pp.Range("D5").Value  = "Immutable"
$pp.Range("D6").Value  = "Only components returned"
$pp.Range("E6").Value  = "Immutable"
$pp.Range("D7").Value  = "Effectively immutable"
$pp.Range("D8").Value  = "Effectively immutable"
$pp.Range("D9").Value  = "Effectively immutable"
$pp.Range("D11").Value = "Effectively immutable"
$pp.Range("D12").Value = "Effectively immutable"
$pp.Range("E12").Value = "Never returned to tenant modules"
$pp.Range("D13").Value = "Never modified"
$pp.Range("E13").Value = "Immutable"
$pp.Range("D14").Value = "Never modified"
$pp.Range("E14").Value = "Effectively immutable"

# widen column E slightly (target stored width ~33.14 chars; COM only
# allows pixel-quantised widths, 32.33 is the closest achievable)
$pp.Columns.Item(5).ColumnWidth = 32.33

# --- "reasons" sheet (sheet5.xml): new reason row ---
$reasons = $wb.Worksheets.Item("reasons")
$reasons.Range("A12").Value = "Only components returned"
$reasons.Range("A12").Select() | Out-Null

# restore "package-private" as the active tab/selection (matches original
# tabSelected + move the cursor to E12 as recorded in the saved view)
$pp.Select() | Out-Null
$pp.Range("E12").Select() | Out-Null
